$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "63×53=" "35×35="
Replace-Text "42×33=" "51×82="
Replace-Text "81×24=" "18×82="
Replace-Text "89×11=" "37×41="
Replace-Text "95×25=" "84×71="
Replace-Text "55×22=" "91×73="
Replace-Text "59×89=" "22×65="
Replace-Text "36×13=" "45×70="
Replace-Text "69×96=" "47×92="
Replace-Text "96×79=" "94×83="
Replace-Text "71×93=" "98×55="
Replace-Text "28×13=" "30×23="
Replace-Text "19×23=" "46×59="
Replace-Text "89×87=" "20×89="
Replace-Text "40×70=" "70×11="
Replace-Text "81×69=" "37×93="
Replace-Text "94×26=" "27×25="
Replace-Text "44×58=" "97×45="
Replace-Text "62×41=" "14×42="
Replace-Text "16×29=" "97×26="
Replace-Text "40×47=" "28×29="
Replace-Text "63×77=" "16×40="
Replace-Text "18×29=" "80×46="
Replace-Text "48×15=" "34×58="
Replace-Text "81×21=" "83×44="
